$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 uses the same style as the other header cells (E1 etc.):
# bold font, centered/top alignment, thin border on all sides.
# Copy the format from E1 (an existing header) so the same style entry is reused.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:52:54.554042",
    "2021-10-05 10:52:54.554052",
    "2021-10-05 10:52:54.554055",
    "2021-10-05 10:52:54.554058",
    "2021-10-05 10:52:54.554061",
    "2021-10-05 10:52:54.554063",
    "2021-10-05 10:52:54.554066",
    "2021-10-05 10:52:54.554068",
    "2021-10-05 10:52:54.554071",
    "2021-10-05 10:52:54.554074",
    "2021-10-05 10:52:54.554077",
    "2021-10-05 10:52:54.554079",
    "2021-10-05 10:52:54.554082",
    "2021-10-05 10:52:54.554084",
    "2021-10-05 10:52:54.554087",
    "2021-10-05 10:52:54.554089",
    "2021-10-05 10:52:54.554092",
    "2021-10-05 10:52:54.554095",
    "2021-10-05 10:52:54.554097",
    "2021-10-05 10:52:54.554100",
    "2021-10-05 10:52:54.554102",
    "2021-10-05 10:52:54.554105",
    "2021-10-05 10:52:54.554108",
    "2021-10-05 10:52:54.554110"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
